$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 65 — appended record.
$ws.Range("A65").Value = "Globo"
$ws.Range("B65").Value = "Bom Dia Inter"
$ws.Range("C65").Value = "Obras"
$ws.Range("D65").Value = "2025-04-07T07:27"
$ws.Range("E65").Value = "Negativo"

# F65 text ends with an embedded newline. Typing it straight into .Value
# (like a live Excel keystroke) would auto-grow the row height; the source
# file's row 65 has no explicit height, so build the string as a formula
# result and paste-special just the value to land a plain literal text
# cell without disturbing row height.
$ws.Range("F65").Formula = '="No encerramento do jornal, imagem de um buraco sinalizado na via próximo à Praça do Amarelinho. *só imagem*."&CHAR(10)'
$ws.Range("F65").Copy()
$ws.Range("F65").PasteSpecial(-4163)
